$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.296.83'
$ws.Range('E2').Value = '  -4.15%  '

$ws.Range('D3').Value = '2.636.24'
$ws.Range('E3').Value = '  -2.34%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('E5').Value = '  -1.38%  '

$ws.Range('E6').Value = '  -1.06%  '

$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('E8').Value = '  -2.02%  '

$ws.Range('E9').Value = '  -0.42%  '

$ws.Range('E10').Value = '  -3.26%  '

$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('E12').Value = '  +1.37%  '

$ws.Range('D13').Value = '3.102.80'
$ws.Range('E13').Value = '  -2.35%  '

$ws.Range('D14').Value = '58.321.64'
$ws.Range('E14').Value = '  -4.04%  '

$ws.Range('E15').Value = '  -2.47%  '

$ws.Range('E16').Value = '  -1.39%  '

$ws.Range('D17').Value = '2.642.48'
$ws.Range('E17').Value = '  -9.10%  '

$ws.Range('E18').Value = '  -3.38%  '

$ws.Range('E19').Value = '  -2.46%  '

$ws.Range('E20').Value = '  -1.64%  '

$ws.Range('E21').Value = '  -2.10%  '

$ws.Range('E23').Value = '  +0.95%  '

$ws.Range('E24').Value = '  +0.63%  '

$ws.Range('E25').Value = '  -1.83%  '

$ws.Range('E26').Value = '  +0.66%  '

$ws.Range('D27').Value = '0.0₃0791'
$ws.Range('E27').Value = '  -2.90%  '

$ws.Range('E28').Value = '  -3.10%  '

$ws.Range('E29').Value = '  -2.90%  '

$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E31').Value = '  +1.43%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E32').Value = '  -1.14%  '

$ws.Range('E33').Value = '  -1.72%  '

$ws.Range('E34').Value = '  -2.83%  '

$ws.Range('E35').Value = '  -4.13%  '

$ws.Range('E36').Value = '  -4.58%  '

$ws.Range('E37').Value = '  -0.71%  '

$ws.Range('E38').Value = '  -2.72%  '

$ws.Range('E39').Value = '  -5.12%  '

$ws.Range('E40').Value = '  -1.19%  '

$ws.Range('E41').Value = '  +0.37%  '

$ws.Range('E42').Value = '  -1.34%  '

$ws.Range('E43').Value = '  -2.33%  '

$ws.Range('E44').Value = '  -5.42%  '

$ws.Range('E45').Value = '  -3.73%  '

$ws.Range('E46').Value = '  +1.60%  '

$ws.Range('E47').Value = '  -0.78%  '

$ws.Range('D48').Value = '2.039.08'
$ws.Range('E48').Value = '  -4.74%  '

$ws.Range('E49').Value = '  -2.46%  '

$ws.Range('E50').Value = '  -3.30%  '

$ws.Range('E51').Value = '  -5.29%  '

# Cells requiring forced text (values that would otherwise be auto-parsed as numbers)
# Force text format, assign, then reset style to Normal to avoid leaving a
# lingering custom-formatted appearance (matches a single shared text xf).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.90'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.15'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.567'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.66'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.72'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '335.98'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.42'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.27'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.07'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.57'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.72'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.58'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.78'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.11'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.17'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.79'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.852'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.603'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0966'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '268.39'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.30'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0534'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.67'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.22'

$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
